$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 148.944201040268
$ws.Range("C2").Value = 3.034028440714535
$ws.Range("D2").Value = 1.369939374923706
$ws.Range("E2").Value = 0.1158927121938705
$ws.Range("B3").Value = 283.8790384769439
$ws.Range("C3").Value = 2.734367804277977
$ws.Range("D3").Value = 1.450635671615601
$ws.Range("E3").Value = 0.1336255080630165
$ws.Range("B4").Value = 559.996348810196
$ws.Range("C4").Value = 3.956669089690096
$ws.Range("D4").Value = 1.46568398475647
$ws.Range("E4").Value = 0.205443587126135
$ws.Range("B5").Value = 147.9574033260346
$ws.Range("C5").Value = 2.581695021344341
$ws.Range("D5").Value = 1.499251413345337
$ws.Range("E5").Value = 0.1597782621201038
$ws.Range("B6").Value = 281.1838232517242
$ws.Range("C6").Value = 2.492547233684098
$ws.Range("D6").Value = 1.407586908340454
$ws.Range("E6").Value = 0.2662355769592323
$ws.Range("B7").Value = 551.9489236831665
$ws.Range("C7").Value = 4.079591150128622
$ws.Range("D7").Value = 1.463583326339722
$ws.Range("E7").Value = 0.07878320146472605
$ws.Range("B8").Value = 144.8851979255676
$ws.Range("C8").Value = 0.3422911666555083
$ws.Range("D8").Value = 1.331832504272461
$ws.Range("E8").Value = 0.0588932411610801
$ws.Range("B9").Value = 288.5357675075531
$ws.Range("C9").Value = 2.455207193314591
$ws.Range("D9").Value = 1.430933856964111
$ws.Range("E9").Value = 0.1580322018010149
$ws.Range("B10").Value = 553.3909600734711
$ws.Range("C10").Value = 5.009806976928636
$ws.Range("D10").Value = 1.296262454986572
$ws.Range("E10").Value = 0.08734082714920767
$ws.Range("B11").Value = 229.5339345932007
$ws.Range("C11").Value = 2.905084350821983
$ws.Range("D11").Value = 1.409769201278686
$ws.Range("E11").Value = 0.1727590792253888
$ws.Range("B12").Value = 444.119971036911
$ws.Range("C12").Value = 1.608215852750877
$ws.Range("D12").Value = 1.566493463516235
$ws.Range("E12").Value = 0.2785014000709363
$ws.Range("B13").Value = 880.7481705665589
$ws.Range("C13").Value = 1.313838911703214
$ws.Range("D13").Value = 1.349617719650269
$ws.Range("E13").Value = 0.1104392556565589
$ws.Range("B14").Value = 226.0651173114776
$ws.Range("C14").Value = 1.609629198597352
$ws.Range("D14").Value = 1.272902202606201
$ws.Range("E14").Value = 0.03110624842315952
$ws.Range("B15").Value = 445.571532535553
$ws.Range("C15").Value = 1.313427412940254
$ws.Range("D15").Value = 1.28677225112915
$ws.Range("E15").Value = 0.03368199035833578
$ws.Range("B16").Value = 885.5678065299987
$ws.Range("C16").Value = 5.325063331153295
$ws.Range("D16").Value = 1.398493957519531
$ws.Range("E16").Value = 0.1520927508007015
$ws.Range("B17").Value = 230.5094635009766
$ws.Range("C17").Value = 1.78845432114904
$ws.Range("D17").Value = 1.302423524856567
$ws.Range("E17").Value = 0.04198255267955443
$ws.Range("B18").Value = 447.0289479255676
$ws.Range("C18").Value = 3.221131304256178
$ws.Range("D18").Value = 1.35668568611145
$ws.Range("E18").Value = 0.1222374296409884
$ws.Range("B19").Value = 719.145511007309
$ws.Range("C19").Value = 24.68404226230715
$ws.Range("D19").Value = 0.8876883029937744
$ws.Range("E19").Value = 0.1238764180127826
